$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the "S:V" block of cell data one column to the right (-> T:W) ---
$src1 = $ws.Range("S1:V184")
$dst1 = $ws.Range("T1:W184")
$src1.Copy($dst1)
$ws.Range("S1:S184").ClearContents()

# --- 2. Shift the "AB:AN" block of cell data two columns to the right (-> AD:AP) ---
$src2 = $ws.Range("AB1:AN184")
$dst2 = $ws.Range("AD1:AP184")
$src2.Copy($dst2)
$ws.Range("AB1:AC184").ClearContents()

# --- 3. Rename shared string "Titulo12" -> "Titulo12_grafico" (now living at AD34) ---
$ws.Range("AD34").Value = "Titulo12_grafico"

# --- 4. Resize the structured tables (ListObjects) to their new locations ---
$ws.ListObjects.Item("Table6").Resize($ws.Range("T2:W4"))
$ws.ListObjects.Item("Table7").Resize($ws.Range("T10:W16"))
$ws.ListObjects.Item("Table8").Resize($ws.Range("T22:U72"))
$ws.ListObjects.Item("Table9").Resize($ws.Range("T78:U128"))
$ws.ListObjects.Item("Table10").Resize($ws.Range("T134:U184"))
$ws.ListObjects.Item("Table11").Resize($ws.Range("AD2:AP13"))
$ws.ListObjects.Item("Table12").Resize($ws.Range("AD19:AP29"))
$ws.ListObjects.Item("Table13").Resize($ws.Range("AD60:AE110"))

# --- 5. Move the second chart picture from column 27 (AB, 0-indexed) to column 29 (AD, 0-indexed) ---
$shapes = $ws.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.Name -eq "Picture 2") {
        $shp.Left = $ws.Range("AD35").Left
    }
}
